$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) <h1>Imię Nazwisko</h1>  ->  <h1>Happyjet Happystyler Tomasz Walczak</h1>
# ------------------------------------------------------------------
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("<h1>Imię Nazwisko</h1>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "<h1>Happyjet Happystyler Tomasz Walczak</h1>"
}

# ------------------------------------------------------------------
# 2) <p class="tagline">Frontend Developer / Grafik / Fotograf</p>
#    -> <p class="tagline">Projektowanie wnętrz</p>
#    (use Range.Text instead of Find's Replace With so straight quotes
#     are not auto-corrected into curly quotes)
# ------------------------------------------------------------------
$rng2 = $d.Content.Duplicate
$found2 = $rng2.Find.Execute('<p class="tagline">Frontend Developer / Grafik / Fotograf</p>', $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Text = '<p class="tagline">Projektowanie wnętrz</p>'
}

# ------------------------------------------------------------------
# 3) Merge the two paragraphs
#      "        Krótki opis kim jesteś, czym się zajmujesz i co robisz najlepiej."
#      "        2–3 zdania, bez lania wody."
#    into a single paragraph:
#      "       Projektuję przestrzenie na jedno skinienie"
# ------------------------------------------------------------------
$p1 = $d.Content.Duplicate
$found3 = $p1.Find.Execute("        Krótki opis kim jesteś, czym się zajmujesz i co robisz najlepiej.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$firstEnd = $p1.End

$p2 = $d.Content.Duplicate
$found4 = $p2.Find.Execute("        2–3 zdania, bez lania wody.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$secondEnd = $p2.End

if ($found3 -and $found4) {
    # Delete the paragraph mark that separates the two paragraphs so they
    # become a single paragraph.
    $markRange = $d.Range($firstEnd, $firstEnd + 1)
    $markRange.Delete()

    # The whole (now merged) paragraph's text, replaced in one shot.
    $whole = $d.Range($p1.Start, $secondEnd - 1)
    $whole.Text = "       Projektuję przestrzenie na jedno skinienie"
}
